$d = $word.ActiveDocument

# Version-control history table: only the first data row (version 2.8.1,
# 8 พฤศจิกายน 2564, "แก้ไข") changes. Each name/suffix lives in its own
# run, so edit each run's text individually (rather than the whole cell)
# to preserve the existing run/formatting split.
# wdReplaceOne (Replace=1) only touches the first match in document
# order, which is exactly this row's "ผู้รับผิดชอบ" / "ผู้ตรวจ" cells.

$r1 = $d.Content
$r1.Find.Execute("วิรัตน์", $true, $false, $false, $false, $false, $true, 1, $false, "ณัฐดนัย", 1) | Out-Null

$r2 = $d.Content
$r2.Find.Execute(" (TL)", $true, $false, $false, $false, $false, $true, 1, $false, " (DM)", 1) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("กิตติพศ ", $true, $false, $false, $false, $false, $true, 1, $false, "วิรัตน์", 1) | Out-Null

$r4 = $d.Content
$r4.Find.Execute("(SP)", $true, $false, $false, $false, $false, $true, 1, $false, " (TL)", 1) | Out-Null
